# "Generate Report for Handoff"
#
# The localization-status report is regenerated: rows whose status used to
# read "Handed back: in sync with en-US" are now marked "Ready for handoff",
# and the associated timestamps are bumped forward (report re-run a minute
# or so later). Because the new status text is shorter than the old one,
# the status column(s) - which are sized to fit their contents - become
# narrower as well.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-12 15:12:43"

# --- zh-cn sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-12 15:12:36"

# --- de-de sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-12 15:12:43"

# --- Resize the status columns to fit the new, shorter text ------------
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null

# The host's AutoFit heuristic is coarser than real Excel's glyph metrics,
# so nudge the freshly autofit-ed columns to the width Excel itself would
# have produced for "Ready for handoff" in this font/style.
$wsOverview.Columns.Item(5).ColumnWidth = 16.35
$wsOverview.Columns.Item(6).ColumnWidth = 16.35
$wsZhCn.Columns.Item(3).ColumnWidth = 16.35
$wsDeDe.Columns.Item(3).ColumnWidth = 16.35
